# Regenerate merged AHB files
# ----------------------------------------------------------------------
# The sheet "AHB-Diff" holds a merged comparison of two AHB versions:
# columns A:J describe the old ("FV2310") version, column K holds a
# "diff" marker, and columns L:U describe the new ("FV2404") version.
# This script:
#   1) renames the header row from the generic "_old"/"_new" suffixes
#      to the concrete version tags "_FV2310"/"_FV2404",
#   2) turns the data range into a native Excel Table ("Table1"), and
#   3) freezes the header row so it stays visible while scrolling.
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row (row 1).
# A1:J1 — "<Name>_old" -> "<Name>_FV2310"
$headersFV2310 = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)
for ($i = 0; $i -lt $headersFV2310.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2310[$i]
}

# K1 ("diff") is unchanged.

# L1:U1 — "<Name>_new" -> "<Name>_FV2404"
$headersFV2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
for ($i = 0; $i -lt $headersFV2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2404[$i]
}

# 2) Convert the used range A1:U62 into a native table named "Table1"
#    (headers already renamed above, so the table columns pick them up).
$tableRange = $ws.Range("A1:U62")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $tableRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# 3) Freeze panes so the header row (row 1) stays visible on scroll.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output "AHB header columns renamed, Table1 created, header row frozen."
